$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear is_parent (column I) flags from rows 3-6 (per diff, I3:I6 values removed)
$ws.Range("I3:I6").ClearContents()

# Row 7: Goblin Coin Bank
$ws.Cells.Item(7,1).Value  = "Goblin Coin Bank"
$ws.Cells.Item(7,2).Value  = "The goblins of Vax, an Island far off the south eastern cost of the Surface plane, have come to our lands to set up banks in your kingdoms. These banks can smelt your gold down to gold bars. Each kingdom you own that has this building may smelt up to 2 trillion (100 gold bars valued at 2 billion each) gold. Each bar adds 0.001 (or 0.1%) to your kingdoms defence for an additional 10% defence bonus."
$ws.Cells.Item(7,3).Value  = 1
$ws.Cells.Item(7,4).Value  = 3
$ws.Cells.Item(7,5).Value  = 0
$ws.Cells.Item(7,6).Value  = 4
$ws.Cells.Item(7,7).Value  = "Building Management"
$ws.Cells.Item(7,8).Value  = 2
$ws.Cells.Item(7,9).Value  = 1

# Row 8: Black Smiths Forge
$ws.Cells.Item(8,1).Value  = "Black Smiths Forge"
$ws.Cells.Item(8,2).Value  = "As you level this skill you will unlock a new building, called Cannoneer Shop to create cannons for your kingdom that can also be deployed out to war. The higher this building level the less iron will be used through out the kingdom (up to 35%, 7% per level). This stacks with any cost reduction for the kingdom."
$ws.Cells.Item(8,3).Value  = 5
$ws.Cells.Item(8,4).Value  = 4
$ws.Cells.Item(8,5).Value  = 0.07
$ws.Cells.Item(8,6).Value  = 4
$ws.Cells.Item(8,7).Value  = "Building Research"
$ws.Cells.Item(8,8).Value  = 4
$ws.Cells.Item(8,9).Value  = 1

# Row 9: Cannoneer Shop
$ws.Cells.Item(9,1).Value  = "Cannoneer Shop"
$ws.Cells.Item(9,2).Value  = "Create cannons for your army, these can be deployed like regular units and move much slower than Trebuchets. These can do devastating damage. When it comes to a kingdoms defence, they can only defend at 65% of their bonus defence against cannons and only 25% for walls against cannons."
$ws.Cells.Item(9,3).Value  = 1
$ws.Cells.Item(9,4).Value  = 5
$ws.Cells.Item(9,5).Value  = 0
$ws.Cells.Item(9,6).Value  = 4
$ws.Cells.Item(9,7).Value  = "Black Smiths Forge"
$ws.Cells.Item(9,8).Value  = 5
$ws.Cells.Item(9,9).Value  = 1

# Update column B width (353 -> 473) per diff
# Note: the COM layer adds a constant offset of 5/6 when persisting the
# "width" attribute, so we compensate by subtracting it here.
$ws.Columns.Item(2).ColumnWidth = 472.16666666666667
